$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new data rows right before the existing row 59, shifting the
# old rows 59-76 down to become rows 63-80 (dimension grows from
# A1:T76 to A1:T80).
$ws.Rows.Item(59).Insert()
$ws.Rows.Item(59).Insert()
$ws.Rows.Item(59).Insert()
$ws.Rows.Item(59).Insert()

# Common columns (A, B, C, E, F, G, H, I, J) are identical for every
# Damasco / Mercado Mayorista Lo Valledor de Santiago row, so copy them
# straight from the row above (row 58) into the 4 freshly inserted rows.
for ($r = 59; $r -le 62; $r++) {
    $ws.Cells.Item($r, 1).Value = 6
    $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100103
    $ws.Cells.Item($r, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($r, 9).Value = 100103003
    $ws.Cells.Item($r, 10).Value = "Damasco"
}

# Row 59
$ws.Cells.Item(59, 4).Value = 44543
$ws.Cells.Item(59, 11).Value = "Castle Brite"
$ws.Cells.Item(59, 12).Value = "Especial"
$ws.Cells.Item(59, 13).Value = 200
$ws.Cells.Item(59, 14).Value = 17000
$ws.Cells.Item(59, 15).Value = 17000
$ws.Cells.Item(59, 16).Value = 17000
$ws.Cells.Item(59, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(59, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(59, 19).Value = 944
$ws.Cells.Item(59, 20).Value = 18

# Row 60
$ws.Cells.Item(60, 4).Value = 44543
$ws.Cells.Item(60, 11).Value = "Castle Brite"
$ws.Cells.Item(60, 12).Value = "Primera"
$ws.Cells.Item(60, 13).Value = 300
$ws.Cells.Item(60, 14).Value = 14000
$ws.Cells.Item(60, 15).Value = 15000
$ws.Cells.Item(60, 16).Value = 14500
$ws.Cells.Item(60, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(60, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(60, 19).Value = 806
$ws.Cells.Item(60, 20).Value = 18

# Row 61
$ws.Cells.Item(61, 4).Value = 44543
$ws.Cells.Item(61, 11).Value = "Castle Brite"
$ws.Cells.Item(61, 12).Value = "Segunda"
$ws.Cells.Item(61, 13).Value = 235
$ws.Cells.Item(61, 14).Value = 10000
$ws.Cells.Item(61, 15).Value = 10000
$ws.Cells.Item(61, 16).Value = 10000
$ws.Cells.Item(61, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(61, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(61, 19).Value = 556
$ws.Cells.Item(61, 20).Value = 18

# Row 62
$ws.Cells.Item(62, 4).Value = 44543
$ws.Cells.Item(62, 11).Value = "Dina"
$ws.Cells.Item(62, 12).Value = "Especial"
$ws.Cells.Item(62, 13).Value = 1100
$ws.Cells.Item(62, 14).Value = 16000
$ws.Cells.Item(62, 15).Value = 16000
$ws.Cells.Item(62, 16).Value = 16000
$ws.Cells.Item(62, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(62, 18).Value = "Región Metropolitana"
$ws.Cells.Item(62, 19).Value = 1000
$ws.Cells.Item(62, 20).Value = 16

Write-Output "applied"
